$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain text formatting (values include dotted
# thousand separators, leading zeros, percent signs and padding spaces that
# must not be reinterpreted as numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '69.007.24'
$ws.Range("E2").Value = '  -5.26%  '

$ws.Range("D3").Value = '3.743.10'
$ws.Range("E3").Value = '  -6.15%  '

$ws.Range("E4").Value = '  +0.42%  '

$ws.Range("D5").Value = '576.79'
$ws.Range("E5").Value = '  -2.64%  '

$ws.Range("D6").Value = '160.79'
$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("D7").Value = '0.646'
$ws.Range("E7").Value = '  -5.44%  '

$ws.Range("E8").Value = '  +0.39%  '

$ws.Range("D9").Value = '0.721'
$ws.Range("E9").Value = '  -3.69%  '

$ws.Range("D10").Value = '0.168'
$ws.Range("E10").Value = '  +0.30%  '

$ws.Range("D11").Value = '51.12'
$ws.Range("E11").Value = '  -4.88%  '

$ws.Range("D12").Value = '0.0000309'
$ws.Range("E12").Value = '  -2.84%  '

$ws.Range("D13").Value = '10.82'
$ws.Range("E13").Value = '  -1.21%  '

$ws.Range("D14").Value = '4.373.26'
$ws.Range("E14").Value = '  -5.35%  '

$ws.Range("D15").Value = '3.793.17'
$ws.Range("E15").Value = '  -5.21%  '

$ws.Range("D16").Value = '20.15'
$ws.Range("E16").Value = '  -1.07%  '

$ws.Range("D17").Value = '1.17'
$ws.Range("E17").Value = '  -8.10%  '

$ws.Range("D18").Value = '13.31'
$ws.Range("E18").Value = '  -5.78%  '

$ws.Range("E19").Value = '  -2.78%  '

$ws.Range("D20").Value = '69.030.09'
$ws.Range("E20").Value = '  -4.83%  '

$ws.Range("D21").Value = '424.85'
$ws.Range("E21").Value = '  -1.90%  '

$ws.Range("D22").Value = '4.57'
$ws.Range("E22").Value = '  -4.60%  '

$ws.Range("D23").Value = '90.81'
$ws.Range("E23").Value = '  -5.60%  '

$ws.Range("D24").Value = '3.17'
$ws.Range("E24").Value = '  -7.67%  '

$ws.Range("D25").Value = '13.42'
$ws.Range("E25").Value = '  -5.27%  '

$ws.Range("D26").Value = '10.88'
$ws.Range("E26").Value = '  -3.65%  '

$ws.Range("D27").Value = '3.83'
$ws.Range("E27").Value = '  -13.24%  '

$ws.Range("E28").Value = '  -0.53%  '

$ws.Range("D29").Value = '10.03'
$ws.Range("E29").Value = '  -4.20%  '

$ws.Range("D30").Value = '33.84'
$ws.Range("E30").Value = '  -6.92%  '

$ws.Range("D31").Value = '7.80'
$ws.Range("E31").Value = '  -1.01%  '

$ws.Range("D32").Value = '13.04'
$ws.Range("E32").Value = '  -5.23%  '

$ws.Range("D33").Value = '46.52'
$ws.Range("E33").Value = '  -4.71%  '

$ws.Range("D34").Value = '0.121'
$ws.Range("E34").Value = '  -7.49%  '

$ws.Range("D35").Value = '67.59'
$ws.Range("E35").Value = '  -4.09%  '

$ws.Range("D36").Value = '0.0₃0944'
$ws.Range("E36").Value = '  +7.78%  '

$ws.Range("D37").Value = '613.85'
$ws.Range("E37").Value = '  -8.44%  '

$ws.Range("D38").Value = '0.411'
$ws.Range("E38").Value = '  -6.03%  '

$ws.Range("E39").Value = '  -0.06%  '

$ws.Range("E40").Value = '  +0.32%  '

$ws.Range("D41").Value = '0.140'
$ws.Range("E41").Value = '  -4.11%  '

$ws.Range("D42").Value = '3.13'
$ws.Range("E42").Value = '  -6.25%  '

$ws.Range("D43").Value = '3.07'
$ws.Range("E43").Value = '  +16.42%  '

$ws.Range("D44").Value = '0.0454'
$ws.Range("E44").Value = '  -7.33%  '

$ws.Range("D45").Value = '2.71'
$ws.Range("E45").Value = '  +4.03%  '

$ws.Range("D46").Value = '9.58'
$ws.Range("E46").Value = '  -10.11%  '

$ws.Range("D47").Value = '0.140'
$ws.Range("E47").Value = '  -6.85%  '

$ws.Range("D48").Value = '2.77'
$ws.Range("E48").Value = '  -17.02%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.781.53'
$ws.Range("E49").Value = '  -2.56%  '

$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").Value = '3.20'
$ws.Range("E50").Value = '  -7.68%  '

$ws.Range("D51").Value = '0.000264'
$ws.Range("E51").Value = '  -1.48%  '
